# Re-normalize the endmember calibration data (columns A/B = Red/NIR)
# against the "White Card" reference row (row 16), and recompute the
# derived NDVI column (C) from the freshly-normalized Red/NIR values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cache the original (pre-edit) reference values from the White Card row
# before any cells in that row get overwritten.
$a16 = $ws.Range("A16").Value2
$b16 = $ws.Range("B16").Value2

for ($r = 2; $r -le 16; $r++) {
    $a = $ws.Range("A$r").Value2
    $b = $ws.Range("B$r").Value2

    $newA = $a / $a16
    $newB = $b / $b16
    $newC = ($newB - $newA) / ($newB + $newA)

    $ws.Range("A$r").Value = $newA
    $ws.Range("B$r").Value = $newB
    $ws.Range("C$r").Value = $newC
}
